# Commit: "Changed bio PARAMETRAR variable info"
# Adds new equipment/method metadata columns (D, E, M, N, O, P, Q) for a
# number of rows in the PARAMETRAR sheet describing sample storage,
# preparation container, weighing method/instrument, etc.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PARAMETRAR")

# --- Rows 110-111 : isotope rows (D13C/D15N UCD) ---
foreach ($r in 110,111) {
    $ws.Range("D$r").Value = "ACES"
    $ws.Range("M$r").Value = "FRYSTO_TENNKAP"
    $ws.Range("N$r").Value = "PLAST"
    $ws.Range("O$r").Value = "UCD_SIA"
    $ws.Range("P$r").Value = "Nej"
    $ws.Range("Q$r").Value = "EA-IRMS"
}

# --- Rows 119-125 : length/weight measurement rows ---
foreach ($r in 119,120,121,122,123,124,125) {
    $ws.Range("D$r").Value = "ACES"
    $ws.Range("E$r").Value = "EJ_REL"
    $ws.Range("M$r").Value = "ACHEXDEE"
    $ws.Range("N$r").Value = "GLAS"
    $ws.Range("O$r").Value = "EJ_REL"
    $ws.Range("P$r").Value = "Nej"
    $ws.Range("Q$r").Value = "VAG"
}

# --- Rows 126-131 : age related rows (NRM) ---
$ws.Range("D126").Value = "NRM"
$ws.Range("E126").Value = "EJ_REL"
$ws.Range("Q126").Value = "LINJAL"

$ws.Range("D127").Value = "NRM"
$ws.Range("E127").Value = "EJ_REL"
$ws.Range("M127").Value = "EJ_REL"
$ws.Range("N127").Value = "EJ_REL"
$ws.Range("O127").Value = "EJ_REL"
$ws.Range("P127").Value = "Vet_ej"
$ws.Range("Q127").Value = "LINJAL"

$ws.Range("D128").Value = "NRM"
$ws.Range("E128").Value = "EJ_REL"
$ws.Range("Q128").Value = "VAG"

$ws.Range("D129").Value = "NRM"
$ws.Range("E129").Value = "EJ_REL"
$ws.Range("M129").Value = "EJ_REL"
$ws.Range("N129").Value = "EJ_REL"
$ws.Range("O129").Value = "EJ_REL"
$ws.Range("P129").Value = "Vet_ej"
$ws.Range("Q129").Value = "VAG"

$ws.Range("D130").Value = "NRM"
$ws.Range("E130").Value = "EJ_REL"
$ws.Range("Q130").Value = "Stereomikroskop"

$ws.Range("D131").Value = "NRM"
$ws.Range("E131").Value = "EJ_REL"
$ws.Range("M131").Value = "EJ_REL"
$ws.Range("N131").Value = "EJ_REL"
$ws.Range("O131").Value = "SAKNAS"
$ws.Range("P131").Value = "Nej"
$ws.Range("Q131").Value = "Stereomikroskop"
